$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper function: set a cell's value as literal TEXT (avoiding Excel's
# automatic number inference for numeric-looking strings like "2.0"),
# while preserving the target cell's existing style/formatting.
# We do this by writing the text (formatted as Text) into a scratch cell,
# copying it, and pasting only the VALUE into the destination.
function Set-TextValue($range, [string]$text) {
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $helper.Clear() | Out-Null
}

# CKD Stage rows (25-27) are being reshuffled: a previously-missing
# "0.0" category row is re-inserted, shifting the existing rows down
# by one position (the old "Absent" label/values are replaced by the
# numeric-looking "0.0" category).

# Row 25 (CKD Stage) - now holds what used to be row 27's data, with
# the "Absent" label replaced by "0.0"
Set-TextValue $ws.Range("B25") "0.0"
$ws.Range("C25").Value = "469 (90.0)"
$ws.Range("D25").Value = "2411 (93.0)"

# Row 26 - now holds what used to be row 25's data
Set-TextValue $ws.Range("B26") "2.0"
$ws.Range("C26").Value = "7 (1.3)"
$ws.Range("D26").Value = "21 (0.8)"

# Row 27 - now holds what used to be row 26's data
Set-TextValue $ws.Range("B27") "3.0"
$ws.Range("C27").Value = "45 (8.6)"
$ws.Range("D27").Value = "161 (6.2)"
